$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ast_config")

# Update the ast_condition status for both data rows from COMPLETE to Queued.
# The leading apostrophe preserves the existing "quote prefix" text style (s="12")
# on these cells instead of Excel re-evaluating/re-styling the value.
$ws.Range("M2").Value = "'Queued"
$ws.Range("M3").Value = "'Queued"
